# Add a new data row (row 36) for 2025-09-20, appending to the
# portfolio-updates time series with the same values as the prior row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so the date-like text is stored as a literal
# string (matching the existing "yyyy-mm-dd" text cells in column A)
# instead of being auto-converted into a date serial number.
$ws.Range("A36").Value = "'2025-09-20"
$ws.Range("B36").Value = 60.40000152587891
$ws.Range("C36").Value = 707.4500122070312
$ws.Range("D36").Value = 336.5499877929688

# Copy the formatting (plain/default style, no quote-prefix marker) from
# the row above so the new cell matches the look of the rest of column A.
$ws.Range("A35").Copy()
$ws.Range("A36").PasteSpecial(-4122)
